$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.878.65"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.497.09"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.582"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("E9").Value = "  +3.76%  "

$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "4.103.44"
$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.61%  "

$ws.Range("D15").Value = "66.885.64"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "3.525.63"
$ws.Range("E17").Value = "  +1.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "389.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("E25").Value = "  -1.05%  "

$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.180"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.90%  "

$ws.Range("E31").Value = "  -1.90%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.876"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.22%  "

$ws.Range("E38").Value = "  -0.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "2.836.96"
$ws.Range("E41").Value = "  +2.01%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0733"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("E46").Value = "  +0.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "339.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.46%  "

$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.841"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
